# Natmi following Dr Hou advice
# Update LR-pair metrics (Arf1-Insr) for rows 2-17 on the active sheet.
# Only columns E, G, H, I, J, K, M, N, O, P, Q, R, S, T change; A-D, F, L stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) whose values are replaced for every data row
$cols = @("E","G","H","I","J","K","M","N","O","P","Q","R","S","T")

# Row => values for each of the columns above (as strings to preserve full precision)
$rowData = @{
  2  = @("3","50.63667066666667","151.910012","0.20081482031288","0.20081482031288","3","13.61877133333333","40.856314","0.3264056993691278","0.3264056993691277","689.609238890641","6206.483150015769","0.06554710186791134","0.06554710186791131")
  3  = @("3","50.63667066666667","151.910012","0.20081482031288","0.20081482031288","3","10.92748633333333","32.782459","0.2619027613928843","0.2619027613928842","553.331526675501","4979.983740079509","0.05259395596855915","0.05259395596855913")
  4  = @("3","50.63667066666667","151.910012","0.20081482031288","0.20081482031288","3","7.492675333333334","22.478026","0.1795794842620271","0.1795794842620271","379.4041332662569","3414.637199396312","0.03606222186395864","0.03606222186395863")
  5  = @("3","50.63667066666667","151.910012","0.20081482031288","0.20081482031288","3","9.684515333333332","29.053546","0.2321120549759609","0.2321120549759609","490.3916135002835","4413.524521502552","0.04661154061245093","0.04661154061245092")
  6  = @("3","67.324","201.972","0.2669934019110801","0.2669934019110801","3","13.61877133333333","40.856314","0.3264056993691278","0.3264056993691277","916.8701612453334","8251.831451208001","0.08714816807772872","0.08714816807772871")
  7  = @("3","67.324","201.972","0.2669934019110801","0.2669934019110801","3","10.92748633333333","32.782459","0.2619027613928843","0.2619027613928842","735.6820899053333","6621.138809148001","0.06992630923419206","0.06992630923419205")
  8  = @("3","67.324","201.972","0.2669934019110801","0.2669934019110801","3","7.492675333333334","22.478026","0.1795794842620271","0.1795794842620271","504.4368741413334","4539.931867272","0.04794653741655589","0.04794653741655588")
  9  = @("3","67.324","201.972","0.2669934019110801","0.2669934019110801","3","9.684515333333332","29.053546","0.2321120549759609","0.2321120549759609","652.0003103013332","5868.002792712","0.06197238718260346","0.06197238718260346")
  10 = @("3","71.72398199999999","215.171946","0.2844428428612245","0.2844428428612245","3","13.61877133333333","40.856314","0.3264056993691278","0.3264056993691277","976.792509974116","8791.132589767045","0.09284376505466088","0.09284376505466087")
  11 = @("3","71.72398199999999","215.171946","0.2844428428612245","0.2844428428612245","3","10.92748633333333","32.782459","0.2619027613928843","0.2619027613928842","783.762833077246","7053.865497695214","0.07449636600379694","0.07449636600379693")
  12 = @("3","71.72398199999999","215.171946","0.2844428428612245","0.2844428428612245","3","7.492675333333334","22.478026","0.1795794842620271","0.1795794842620271","537.404510739844","4836.640596658596","0.05108009902304351","0.05108009902304349")
  13 = @("3","71.72398199999999","215.171946","0.2844428428612245","0.2844428428612245","3","9.684515333333332","29.053546","0.2321120549759609","0.2321120549759609","694.6120034467239","6251.508031020515","0.06602261277972314","0.06602261277972314")
  14 = @("3","62.47139133333334","187.414174","0.2477489349148154","0.2477489349148154","3","13.61877133333333","40.856314","0.3264056993691278","0.3264056993691277","850.7835934438485","7657.052340994637","0.08086666436882686","0.08086666436882683")
  15 = @("3","62.47139133333334","187.414174","0.2477489349148154","0.2477489349148154","3","10.92748633333333","32.782459","0.2619027613928843","0.2619027613928842","682.6552750193185","6143.897475173867","0.06488613018633613","0.0648861301863361")
  16 = @("3","62.47139133333334","187.414174","0.2477489349148154","0.2477489349148154","3","7.492675333333334","22.478026","0.1795794842620271","0.1795794842620271","468.0778528822805","4212.700675940524","0.04449062595846908","0.04449062595846906")
  17 = @("3","62.47139133333334","187.414174","0.2477489349148154","0.2477489349148154","3","9.684515333333332","29.053546","0.2321120549759609","0.2321120549759609","605.0051472623337","5445.046325361003","0.05750551440118341","0.0575055144011834")
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $ws.Range("$col$row").Value = [double]$values[$i]
    }
}
